$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H39").Value = 909.4666999999999
$ws.Range("I39").Value = 734.8570999999999
$ws.Range("J39").Value = 1062.25
$ws.Range("K39").Value = 2204.5713
$ws.Range("L39").Value = 3186.75
$ws.Range("M39").Value = -1908.5713
$ws.Range("N39").Value = -3778.75
$ws.Range("H41").Value = 23811212
$ws.Range("I41").Value = 47620630
$ws.Range("J41").Value = 1790
$ws.Range("K41").Value = 47620630
$ws.Range("L41").Value = 1790
$ws.Range("M41").Value = -47620190
$ws.Range("N41").Value = -2670
$ws.Range("H53").Value = 1847.5
$ws.Range("J53").Value = 1778.1666
$ws.Range("L53").Value = 1778.1666
$ws.Range("N53").Value = -3052.1666
$ws.Range("H113").Value = 2616
$ws.Range("I113").Value = 1924
$ws.Range("K113").Value = 1924
$ws.Range("M113").Value = 1330
$ws.Range("H137").Value = 1671
$ws.Range("I137").Value = 2332.3333
$ws.Range("K137").Value = 6996.999899999999
$ws.Range("M137").Value = -4446.999899999999
$ws.Range("H138").Value = 726692.1
$ws.Range("I138").Value = 1154.8235
$ws.Range("J138").Value = 1013532.44
$ws.Range("K138").Value = 3464.4705
$ws.Range("L138").Value = 3040597.32
$ws.Range("M138").Value = 1675.5295
$ws.Range("N138").Value = -3050877.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1523.04
$ws.Range("I61").Value = 1384.5714
$ws.Range("J61").Value = 2250
$ws.Range("K61").Value = 1384.5714
$ws.Range("L61").Value = 2250
$ws.Range("M61").Value = -1172.5714
$ws.Range("N61").Value = -2674
$ws.Range("H74").Value = 1105.1428
$ws.Range("I74").Value = 1105.1428
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1105.1428
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -231.1428000000001
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1105.1428
$ws.Range("I77").Value = 1105.1428
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5525.714
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1157.714
$ws.Range("N77").ClearContents()
$ws.Range("H122").Value = 1222.4
$ws.Range("I122").Value = 1222.4
$ws.Range("K122").Value = 3667.2
$ws.Range("M122").Value = -1217.2
$ws.Range("H136").Value = 1523.04
$ws.Range("I136").Value = 1384.5714
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 4153.7142
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -1603.7142
$ws.Range("N136").Value = -11850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6180.857
$ws.Range("I134").Value = 1193.6666
$ws.Range("J134").Value = 18648.834
$ws.Range("K134").Value = 3580.9998
$ws.Range("L134").Value = 55946.50199999999
$ws.Range("M134").Value = -1045.9998
$ws.Range("N134").Value = -61016.50199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1748.0834
$ws.Range("I31").Value = 1748.0834
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1748.0834
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1453.0834
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1748.0834
$ws.Range("I34").Value = 1748.0834
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1748.0834
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1546.0834
$ws.Range("N34").ClearContents()
$ws.Range("H58").Value = 644.4583
$ws.Range("I58").Value = 645.15
$ws.Range("J58").Value = 641
$ws.Range("K58").Value = 645.15
$ws.Range("L58").Value = 641
$ws.Range("M58").Value = -442.15
$ws.Range("N58").Value = -1047
$ws.Range("H99").Value = 1522.9375
$ws.Range("I99").Value = 1520.7693
$ws.Range("K99").Value = 1520.7693
$ws.Range("M99").Value = -22.76929999999993
$ws.Range("H126").Value = 1522.9375
$ws.Range("I126").Value = 1520.7693
$ws.Range("K126").Value = 4562.3079
$ws.Range("M126").Value = -2092.3079
$ws.Range("H132").Value = 4764.129
$ws.Range("I132").Value = 4959.1924
$ws.Range("J132").Value = 3749.8
$ws.Range("K132").Value = 14877.5772
$ws.Range("L132").Value = 11249.4
$ws.Range("M132").Value = -12347.5772
$ws.Range("N132").Value = -16309.4
$ws.Range("H134").Value = 33335908
$ws.Range("I134").Value = 55558264
$ws.Range("K134").Value = 166674792
$ws.Range("M134").Value = -166672257
$ws.Range("H136").Value = 644.4583
$ws.Range("I136").Value = 645.15
$ws.Range("J136").Value = 641
$ws.Range("K136").Value = 1935.45
$ws.Range("L136").Value = 1923
$ws.Range("M136").Value = 614.5500000000002
$ws.Range("N136").Value = -7023

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 230
$ws.Range("I21").Value = 230
$ws.Range("K21").Value = 690
$ws.Range("M21").Value = -517
$ws.Range("H99").Value = 3016.8
$ws.Range("J99").Value = 3016.8
$ws.Range("L99").Value = 9050.400000000001
$ws.Range("N99").Value = -13542.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6129.75
$ws.Range("I43").Value = 6000
$ws.Range("J43").Value = 6173
$ws.Range("K43").Value = 6000
$ws.Range("L43").Value = 6173
$ws.Range("M43").Value = -5849
$ws.Range("N43").Value = -6475
$ws.Range("H86").Value = 27497
$ws.Range("J86").Value = 27497
$ws.Range("L86").Value = 27497
$ws.Range("N86").Value = -29869
$ws.Range("H89").Value = 27497
$ws.Range("J89").Value = 27497
$ws.Range("L89").Value = 82491
$ws.Range("N89").Value = -94347
$ws.Range("H122").Value = 1173.75
$ws.Range("J122").Value = 1500
$ws.Range("L122").Value = 4500
$ws.Range("N122").Value = -9400
$ws.Range("H126").Value = 3200.95
$ws.Range("I126").Value = 1866.875
$ws.Range("J126").Value = 4090.3333
$ws.Range("K126").Value = 5600.625
$ws.Range("L126").Value = 12270.9999
$ws.Range("M126").Value = -3130.625
$ws.Range("N126").Value = -17210.9999
$ws.Range("H132").Value = 2276.6316
$ws.Range("I132").Value = 2204.7856
$ws.Range("K132").Value = 6614.3568
$ws.Range("M132").Value = -4084.3568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 900
$ws.Range("I4").Value = 900
$ws.Range("K4").Value = 900
$ws.Range("M4").Value = -787
$ws.Range("H7").Value = 1622.7778
$ws.Range("I7").Value = 1296.2858
$ws.Range("K7").Value = 1296.2858
$ws.Range("M7").Value = -1184.2858
$ws.Range("H22").Value = 950.9091
$ws.Range("I22").Value = 949.3333
$ws.Range("J22").Value = 954.2857
$ws.Range("K22").Value = 949.3333
$ws.Range("L22").Value = 954.2857
$ws.Range("M22").Value = -654.3333
$ws.Range("N22").Value = -1544.2857
$ws.Range("H27").Value = 950.9091
$ws.Range("I27").Value = 949.3333
$ws.Range("J27").Value = 954.2857
$ws.Range("K27").Value = 949.3333
$ws.Range("L27").Value = 954.2857
$ws.Range("M27").Value = -842.3333
$ws.Range("N27").Value = -1168.2857
$ws.Range("H28").Value = 900
$ws.Range("I28").Value = 900
$ws.Range("K28").Value = 900
$ws.Range("M28").Value = -668
$ws.Range("H37").Value = 900
$ws.Range("I37").Value = 900
$ws.Range("K37").Value = 900
$ws.Range("M37").Value = -793
$ws.Range("H40").Value = 2912.111
$ws.Range("J40").Value = 3401.6667
$ws.Range("L40").Value = 3401.6667
$ws.Range("N40").Value = -3673.6667
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H126").Value = 1622.7778
$ws.Range("I126").Value = 1296.2858
$ws.Range("K126").Value = 3888.8574
$ws.Range("M126").Value = -1418.8574
$ws.Range("H132").Value = 74378
$ws.Range("J132").Value = 113889.22
$ws.Range("L132").Value = 341667.66
$ws.Range("N132").Value = -346727.66
$ws.Range("H136").Value = 6552.8
$ws.Range("I136").Value = 9028.538
$ws.Range("J136").Value = 1955
$ws.Range("K136").Value = 27085.614
$ws.Range("L136").Value = 5865
$ws.Range("M136").Value = -24535.614
$ws.Range("N136").Value = -10965

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7030543
$ws.Range("I122").Value = 8128543
$ws.Range("J122").Value = 3343
$ws.Range("K122").Value = 24385629
$ws.Range("L122").Value = 10029
$ws.Range("M122").Value = -24383179
$ws.Range("N122").Value = -14929
$ws.Range("H126").Value = 62501412
$ws.Range("I126").Value = 66667976
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 200003928
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -200001458
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 2667.7646
$ws.Range("I132").Value = 2491
$ws.Range("J132").Value = 3037.3635
$ws.Range("K132").Value = 7473
$ws.Range("L132").Value = 9112.0905
$ws.Range("M132").Value = -4943
$ws.Range("N132").Value = -14172.0905
$ws.Range("H136").Value = 531.7143
$ws.Range("I136").Value = 270.33334
$ws.Range("K136").Value = 811.0000200000001
$ws.Range("M136").Value = 1738.99998
